{"js": "// Update the division equations in the table to the new set of values.\n// Each change below is a simple exact-text replacement: the old equation\n// text (e.g. \"495\u00f74=\") is found in the document and replaced with the\n// new equation text (e.g. \"585\u00f78=\").\nconst replacements = [\n  [\"495\u00f74=\", \"585\u00f78=\"],\n  [\"735\u00f74=\", \"978\u00f79=\"],\n  [\"572\u00f73=\", \"655\u00f74=\"],\n  [\"162\u00f78=\", \"532\u00f72=\"],\n  [\"747\u00f73=\", \"415\u00f76=\"],\n  [\"433\u00f73=\", \"171\u00f79=\"],\n  [\"858\u00f76=\", \"339\u00f76=\"],\n  [\"643\u00f78=\", \"794\u00f73=\"],\n  [\"676\u00f72=\", \"540\u00f73=\"],\n  [\"783\u00f76=\", \"717\u00f75=\"],\n  [\"274\u00f76=\", \"133\u00f74=\"],\n  [\"496\u00f79=\", \"852\u00f76=\"],\n  [\"744\u00f72=\", \"702\u00f73=\"],\n  [\"644\u00f72=\", \"204\u00f78=\"],\n  [\"810\u00f76=\", \"955\u00f75=\"],\n  [\"149\u00f77=\", \"837\u00f76=\"],\n  [\"889\u00f79=\", \"576\u00f79=\"],\n  [\"348\u00f74=\", \"732\u00f72=\"],\n  [\"655\u00f76=\", \"826\u00f74=\"],\n  [\"774\u00f74=\", \"436\u00f76=\"],\n  [\"674\u00f79=\", \"988\u00f74=\"],\n  [\"803\u00f74=\", \"634\u00f77=\"],\n  [\"875\u00f79=\", \"603\u00f77=\"],\n  [\"345\u00f78=\", \"238\u00f73=\"],\n  [\"482\u00f78=\", \"160\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division equations in the table to the new set of values.\n# Each change below is a simple exact-text replacement: the old equation\n# text (e.g. \"495\u00f74=\") is found in the document and replaced with the\n# new equation text (e.g. \"585\u00f78=\").\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"495\u00f74=\", \"585\u00f78=\"),\n  @(\"735\u00f74=\", \"978\u00f79=\"),\n  @(\"572\u00f73=\", \"655\u00f74=\"),\n  @(\"162\u00f78=\", \"532\u00f72=\"),\n  @(\"747\u00f73=\", \"415\u00f76=\"),\n  @(\"433\u00f73=\", \"171\u00f79=\"),\n  @(\"858\u00f76=\", \"339\u00f76=\"),\n  @(\"643\u00f78=\", \"794\u00f73=\"),\n  @(\"676\u00f72=\", \"540\u00f73=\"),\n  @(\"783\u00f76=\", \"717\u00f75=\"),\n  @(\"274\u00f76=\", \"133\u00f74=\"),\n  @(\"496\u00f79=\", \"852\u00f76=\"),\n  @(\"744\u00f72=\", \"702\u00f73=\"),\n  @(\"644\u00f72=\", \"204\u00f78=\"),\n  @(\"810\u00f76=\", \"955\u00f75=\"),\n  @(\"149\u00f77=\", \"837\u00f76=\"),\n  @(\"889\u00f79=\", \"576\u00f79=\"),\n  @(\"348\u00f74=\", \"732\u00f72=\"),\n  @(\"655\u00f76=\", \"826\u00f74=\"),\n  @(\"774\u00f74=\", \"436\u00f76=\"),\n  @(\"674\u00f79=\", \"988\u00f74=\"),\n  @(\"803\u00f74=\", \"634\u00f77=\"),\n  @(\"875\u00f79=\", \"603\u00f77=\"),\n  @(\"345\u00f78=\", \"238\u00f73=\"),\n  @(\"482\u00f78=\", \"160\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
